# Apply weekly update: insert 3 new rows (Especial/Primera/Segunda) for the
# latest reporting date (2022-01-24, serial 44585) at the top of the data
# block (row 806), pushing the existing rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 806 (first data row of the
# block being updated). Formatting is inherited from the row below, which
# already carries the date style used throughout the column.
$ws.Rows("806:808").Insert()

# Common / constant values shared by every record in this data set.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100101
$producto   = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "$/bandeja 7 kilos"
$origen      = "Provincia de Melipilla"
$kgUnidad    = 7

$fecha = 44585

# Row 806: Especial
$ws.Range("A806").Value2 = $mercadoId
$ws.Range("B806").Value2 = $mercado
$ws.Range("C806").Value2 = $region
$ws.Range("D806").Value2 = $fecha
$ws.Range("E806").Value2 = $codreg
$ws.Range("F806").Value2 = $tipo
$ws.Range("G806").Value2 = $productoId
$ws.Range("H806").Value2 = $producto
$ws.Range("I806").Value2 = $categoriaId
$ws.Range("J806").Value2 = $categoria
$ws.Range("K806").Value2 = $variedad
$ws.Range("L806").Value2 = "Especial"
$ws.Range("M806").Value2 = 1200
$ws.Range("N806").Value2 = 6000
$ws.Range("O806").Value2 = 7000
$ws.Range("P806").Value2 = 6500
$ws.Range("Q806").Value2 = $unidad
$ws.Range("R806").Value2 = $origen
$ws.Range("S806").Value2 = 929
$ws.Range("T806").Value2 = $kgUnidad

# Row 807: Primera
$ws.Range("A807").Value2 = $mercadoId
$ws.Range("B807").Value2 = $mercado
$ws.Range("C807").Value2 = $region
$ws.Range("D807").Value2 = $fecha
$ws.Range("E807").Value2 = $codreg
$ws.Range("F807").Value2 = $tipo
$ws.Range("G807").Value2 = $productoId
$ws.Range("H807").Value2 = $producto
$ws.Range("I807").Value2 = $categoriaId
$ws.Range("J807").Value2 = $categoria
$ws.Range("K807").Value2 = $variedad
$ws.Range("L807").Value2 = "Primera"
$ws.Range("M807").Value2 = 750
$ws.Range("N807").Value2 = 4000
$ws.Range("O807").Value2 = 5000
$ws.Range("P807").Value2 = 4500
$ws.Range("Q807").Value2 = $unidad
$ws.Range("R807").Value2 = $origen
$ws.Range("S807").Value2 = 643
$ws.Range("T807").Value2 = $kgUnidad

# Row 808: Segunda
$ws.Range("A808").Value2 = $mercadoId
$ws.Range("B808").Value2 = $mercado
$ws.Range("C808").Value2 = $region
$ws.Range("D808").Value2 = $fecha
$ws.Range("E808").Value2 = $codreg
$ws.Range("F808").Value2 = $tipo
$ws.Range("G808").Value2 = $productoId
$ws.Range("H808").Value2 = $producto
$ws.Range("I808").Value2 = $categoriaId
$ws.Range("J808").Value2 = $categoria
$ws.Range("K808").Value2 = $variedad
$ws.Range("L808").Value2 = "Segunda"
$ws.Range("M808").Value2 = 500
$ws.Range("N808").Value2 = 3000
$ws.Range("O808").Value2 = 3500
$ws.Range("P808").Value2 = 3250
$ws.Range("Q808").Value2 = $unidad
$ws.Range("R808").Value2 = $origen
$ws.Range("S808").Value2 = 464
$ws.Range("T808").Value2 = $kgUnidad
